# Update the public EPEX Spot prices workbook:
#  - "Prix Spot" sheet: add a new day column AW (01-aug) with its 24 hourly values
#  - "Gaz" sheet: append a new row (2025-07-30 / 34.175)
#  - "CO2" sheet: append a new row (2025-07-30 / 72.12)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" -- add column AW
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Header cell, reusing the same bold/centered/bordered style as the other
# header cells on row 1 (e.g. AV1) by copying formats from the neighbouring
# header cell, then writing the new label.
$wsSpot.Range("AV1").Copy() | Out-Null
$wsSpot.Range("AW1").PasteSpecial(-4122) | Out-Null
$wsSpot.Range("AW1").Value = "01-aug"

$spotValues = @(
    102.22,
    93.47,
    81,
    72.55,
    63.98,
    73.59999999999999,
    83.14,
    100.62,
    100,
    65.06,
    45.81,
    35.08,
    31.87,
    25.06,
    20.01,
    13.43,
    17.42,
    27.47,
    44.8,
    70.12,
    87.34,
    103.02,
    102.52,
    95.41
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 49).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" -- append row 46
# Force the date column to stay plain text (matching every other row in the
# column, which are inline/shared strings, not real dates) by formatting the
# cell as Text before assigning, then resetting the style so no stray
# number-format sticks to the cell.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A46").NumberFormat = "@"
$wsGaz.Range("A46").Value = "2025-07-30"
$wsGaz.Range("A46").Style = "Normal"
$wsGaz.Range("B46").Value = 34.175

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" -- append row 46
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A46").NumberFormat = "@"
$wsCo2.Range("A46").Value = "2025-07-30"
$wsCo2.Range("A46").Style = "Normal"
$wsCo2.Range("B46").Value = 72.12
